$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 883, shifting existing rows 883:961 down to 884:962
$ws.Rows.Item(883).Insert()

# Populate the new row 883 with values (same categorical columns as surrounding rows)
$ws.Cells.Item(883, 1).Value = 5
$ws.Cells.Item(883, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(883, 3).Value = "Maule"
$ws.Cells.Item(883, 4).Value = 44578
$ws.Cells.Item(883, 5).Value = 7
$ws.Cells.Item(883, 6).Value = "Fruta"
$ws.Cells.Item(883, 7).Value = 100102
$ws.Cells.Item(883, 8).Value = "Cítricos"
$ws.Cells.Item(883, 9).Value = 100102003
$ws.Cells.Item(883, 10).Value = "Limón"
$ws.Cells.Item(883, 11).Value = "Sin especificar"
$ws.Cells.Item(883, 12).Value = "1a plateado"
$ws.Cells.Item(883, 13).Value = 350
$ws.Cells.Item(883, 14).Value = 15000
$ws.Cells.Item(883, 15).Value = 15000
$ws.Cells.Item(883, 16).Value = 15000
$ws.Cells.Item(883, 17).Value = "$/malla 14 kilos"
$ws.Cells.Item(883, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(883, 19).Value = 1071
$ws.Cells.Item(883, 20).Value = 14
